$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("files")

# Assign VS (row 23) Status = InProgress, assigned to Jessica
$ws.Range("C23").Value = "InProgress"
$ws.Range("D23").Value = "Jessica"

$excel.ActiveWindow.ScrollRow = 11
$ws.Range("C24").Select() | Out-Null
